$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 331-332 (everything from old row 331 down shifts to 333-349)
$ws.Rows("331:332").Insert()

# New row 331: Alcachofa / Española / Extra, fecha 2021-11-09 (serial 44509), Región Metropolitana
$ws.Cells.Item(331, 1).Value = 9
$ws.Cells.Item(331, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(331, 3).Value = "Metropolitana"
$ws.Cells.Item(331, 4).Value = 44509
$ws.Cells.Item(331, 5).Value = 13
$ws.Cells.Item(331, 6).Value = 100112013
$ws.Cells.Item(331, 7).Value = "Alcachofa"
$ws.Cells.Item(331, 8).Value = "Española"
$ws.Cells.Item(331, 9).Value = "Extra"
$ws.Cells.Item(331, 10).Value = 61
$ws.Cells.Item(331, 11).Value = 10000
$ws.Cells.Item(331, 12).Value = 11000
$ws.Cells.Item(331, 13).Value = 10492
$ws.Cells.Item(331, 14).Value = "$/caja 25 unidades"
$ws.Cells.Item(331, 15).Value = "Región Metropolitana"
$ws.Cells.Item(331, 16).Value = 10492
$ws.Cells.Item(331, 17).Value = 1
$ws.Cells.Item(331, 18).Value = "Hortaliza"

# New row 332: Alcachofa / Española / Primera, fecha 2021-11-09 (serial 44509), Región Metropolitana
$ws.Cells.Item(332, 1).Value = 9
$ws.Cells.Item(332, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(332, 3).Value = "Metropolitana"
$ws.Cells.Item(332, 4).Value = 44509
$ws.Cells.Item(332, 5).Value = 13
$ws.Cells.Item(332, 6).Value = 100112013
$ws.Cells.Item(332, 7).Value = "Alcachofa"
$ws.Cells.Item(332, 8).Value = "Española"
$ws.Cells.Item(332, 9).Value = "Primera"
$ws.Cells.Item(332, 10).Value = 79
$ws.Cells.Item(332, 11).Value = 8000
$ws.Cells.Item(332, 12).Value = 9000
$ws.Cells.Item(332, 13).Value = 8506
$ws.Cells.Item(332, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(332, 15).Value = "Región Metropolitana"
$ws.Cells.Item(332, 16).Value = 284
$ws.Cells.Item(332, 17).Value = 30
$ws.Cells.Item(332, 18).Value = "Hortaliza"
